$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 112379125
$ws.Range("B3").Value = 90812
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 4366
$ws.Range("F3").Value = "Skarp dropptaggsvamp"
$ws.Range("G3").Value = "Hydnellum peckii"
$ws.Range("H3").Value = "Banker"
$ws.Range("P3").Value = "Lortmossen, Vstm"
$ws.Range("Q3").Value = 531963
$ws.Range("R3").Value = 6622561
$ws.Range("S3").Value = 25
$ws.Range("T3").Value = "Västmanland"
$ws.Range("U3").Value = "Skinnskatteberg"
$ws.Range("V3").Value = "Västmanland"
$ws.Range("W3").Value = "Skinnskatteberg"
$ws.Range("Y3").Value = "'2023-08-27"
$ws.Range("AA3").Value = "'2023-09-18"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AW3").Value = "Mikael Hagström"
$ws.Range("AX3").Value = "Mikael Hagström"

# Row 4
$ws.Range("A4").Value = 112379126
$ws.Range("B4").Value = 90812
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 4366
$ws.Range("F4").Value = "Skarp dropptaggsvamp"
$ws.Range("G4").Value = "Hydnellum peckii"
$ws.Range("H4").Value = "Banker"
$ws.Range("P4").Value = "Lortmossen, Vstm"
$ws.Range("Q4").Value = 532036
$ws.Range("R4").Value = 6622648
$ws.Range("S4").Value = 25
$ws.Range("T4").Value = "Västmanland"
$ws.Range("U4").Value = "Skinnskatteberg"
$ws.Range("V4").Value = "Västmanland"
$ws.Range("W4").Value = "Skinnskatteberg"
$ws.Range("Y4").Value = "'2023-08-27"
$ws.Range("AA4").Value = "'2023-09-18"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = "Mikael Hagström"
$ws.Range("AX4").Value = "Mikael Hagström"
